# "Update Data Sources from LFX" — the automated data refresh re-created the
# tables on the data-source slides, which reset each table's style from the
# deck's custom "Table_0" style back to PowerPoint's built-in default table
# style ({19266D9C-3169-4CB3-B0A4-9DD695BD7118}).
#
# Walk every slide/shape in the deck and, for every shape that is a table,
# re-apply the built-in style so the <a:tableStyleId> in the OOXML matches.

$p = $ppt.ActivePresentation
$newTableStyleId = "{19266D9C-3169-4CB3-B0A4-9DD695BD7118}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
